# "pasted initial requirements for traceability"
# 1) Rename Sheet1 -> "Competency Questions"
# 2) Add a new sheet "Initial Requirements" after it, with an ID/Requirement/
#    refers-to/source traceability table.

$wb = $excel.ActiveWorkbook

$cq = $wb.Worksheets.Item(1)
$cq.Name = "Competency Questions"

$req = $wb.Worksheets.Add($null, $cq)
$req.Name = "Initial Requirements"

# --- values, entered in the same order the original author typed/pasted
# them (column A IDs first for R1-R4, then the whole Requirement column,
# then the remaining IDs R5-R7, then the "refers to" column, then "source") ---

$req.Cells.Item(1, 1).Value = "ID"

$req.Cells.Item(1, 2).Value = "Requirement"
$req.Cells.Item(2, 1).Value = "R1"
$req.Cells.Item(3, 1).Value = "R2"
$req.Cells.Item(4, 1).Value = "R3"
$req.Cells.Item(5, 1).Value = "R4"

$req.Cells.Item(2, 2).Value = "User should be capable of performing every aspect in the normative GSN v3 standard with OntoGSN."
$req.Cells.Item(3, 2).Value = "OntoGSN should represent knowledge about GSN assurance cases in direct conformity with the standard metamodel."
$req.Cells.Item(4, 2).Value = "OntoGSN should be aligned with the ASCE XML schema."
$req.Cells.Item(5, 2).Value = "The user should be able to use OntoGSN as a standalone framework."
$req.Cells.Item(6, 2).Value = "The user should be able to use OntoGSN in an accessible manner."
$req.Cells.Item(7, 2).Value = "OntoGSN should not duplicate the features provided by existing tools and frameworks (e.g., GUI, imports/exports)."
$req.Cells.Item(8, 2).Value = "The user should have an accessible and easy way to create SPARQL queries based on a case ABox conformant to the OntoGSN TBox."

$req.Cells.Item(6, 1).Value = "R5"
$req.Cells.Item(7, 1).Value = "R6"
$req.Cells.Item(8, 1).Value = "R7"

$req.Cells.Item(1, 3).Value = "refers to"
$req.Cells.Item(2, 3).Value = "ontology"
$req.Cells.Item(3, 3).Value = "ontology"
$req.Cells.Item(4, 3).Value = "ontology"
$req.Cells.Item(5, 3).Value = "framework, supporting tools"
$req.Cells.Item(6, 3).Value = "framework, supporting tools"
$req.Cells.Item(7, 3).Value = "supporting tools"
$req.Cells.Item(8, 3).Value = "supporting tools"

$req.Cells.Item(1, 4).Value = "source"
$req.Cells.Item(2, 4).Value = "GSN Standard"
$req.Cells.Item(3, 4).Value = "GSN Standard"
$req.Cells.Item(7, 4).Value = "https://scsc.uk/index.php/tools"
$req.Cells.Item(4, 4).Value = "https://scsc.uk/index.php/tools; feedback"
$req.Cells.Item(6, 4).Value = "feedback"

# --- formatting ---

$req.Range("A1:D1").Font.Bold = $true
$req.Range("B2:B8").WrapText = $true

# paste-introduced "no fill" marker on D4 (matches the source workbook's
# extra cellXfs entry with applyFill but no color change)
$req.Cells.Item(4, 4).Interior.ColorIndex = -4142

$req.Columns.Item(2).ColumnWidth = 38
$req.Columns.Item(3).ColumnWidth = 24.6

# --- view state: Initial Requirements becomes the active/selected tab ---

$cq.Range("C19").Select() | Out-Null
$req.Activate() | Out-Null
$req.Range("C11").Select() | Out-Null
